$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 ("land"): drop the first data row (index 14 / 南投縣仁愛郷廬山段
# 06760000地號), then add two new trailing columns: "portion" (numeric
# share, e.g. 2分之1 -> 0.5, 全部 -> 1) and "total" (area * portion).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Remove the obsolete first land record entirely - remaining rows shift up.
$ws1.Rows.Item(2).Delete()

# Figure out how many data rows remain.
$lastRow = $ws1.UsedRange.Rows.Count

# New header cells, copying the header formatting (bold/border/center) from
# the existing last header cell so the new columns look consistent.
$ws1.Range("O1").Copy()
$ws1.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws1.Cells.Item(1, 16).Value = "portion"
$ws1.Cells.Item(1, 17).Value = "total"

for ($r = 2; $r -le $lastRow; $r++) {
    $area = $ws1.Cells.Item($r, 3).Value()
    $shareText = $ws1.Cells.Item($r, 4).Value()

    if ($shareText -eq "2分之1") {
        $portion = 0.5
    } else {
        $portion = 1
    }

    $ws1.Cells.Item($r, 16).Value = $portion
    $ws1.Cells.Item($r, 17).Value = $area * $portion
}

# ---------------------------------------------------------------------------
# Sheet 2 ("car"): the sheet had a header row (row 1) plus a single data row
# (row 2, with an index value in column A). Fold the data up into row 1
# (keeping row 1's header formatting), drop column A's index value, and
# remove the now-empty former data row.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

for ($c = 2; $c -le 7; $c++) {
    $ws2.Cells.Item(1, $c).Value = $ws2.Cells.Item(2, $c).Value()
}
$ws2.Rows.Item(2).Delete()
$ws2.Range("A1").ClearContents()
